$wb = $excel.ActiveWorkbook

# --- Sheet management -------------------------------------------------
# Rename the existing "Submissions" sheet to "cv-data" and add a new,
# empty "test" sheet right after it.
$ws1 = $wb.ActiveSheet
$ws1.Name = "cv-data"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "test"

# --- Trim "cv-data" down to the new shape (A1:D2) ----------------------
# Original sheet had columns id/name/email/phone/message/submittedAt and
# two data rows. The new layout only needs id/submittedAt/Name/Email and
# a single data row, so drop the extra columns (E:F) and the extra row (3).
$ws1.Range("E1:F3").Delete()
$ws1.Rows.Item(3).Delete()

# --- Update the header row ---------------------------------------------
$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "submittedAt"
$ws1.Range("C1").Value = "Name"
$ws1.Range("D1").Value = "Email"

# --- Update the single remaining data row -------------------------------
$ws1.Range("A2").Value = "8eaf02eb-5f97-437e-9176-9c1f4dad580e"
$ws1.Range("B2").Value = "2025-04-26T14:47:40.065Z"
$ws1.Range("C2").Value = "minal"
$ws1.Range("D2").Value = "minal@gmail.com"
